$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 5931
$ws.Range("F8").Value = 10109
$ws.Range("F12").Value = 4002
$ws.Range("F18").Value = 687
$ws.Range("F19").Value = 3986
$ws.Range("F22").Value = 5578
$ws.Range("F27").Value = 8341
$ws.Range("F31").Value = 2270
$ws.Range("F34").Value = 1857
$ws.Range("F40").Value = 18
$ws.Range("F43").Value = 47
$ws.Range("F44").Value = 70
$ws.Range("F45").Value = 195
$ws.Range("F46").Value = 1385
$ws.Range("F47").Value = 248
$ws.Range("F49").Value = 12

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 5931
$ws.Range("F9").Value = 10109
$ws.Range("F11").Value = 4002
$ws.Range("F18").Value = 687
$ws.Range("F19").Value = 3986
$ws.Range("F22").Value = 5578
$ws.Range("F27").Value = 8341
$ws.Range("F31").Value = 2270
$ws.Range("F34").Value = 1857
$ws.Range("F39").Value = 18
$ws.Range("F42").Value = 47
$ws.Range("F43").Value = 70
$ws.Range("F44").Value = 195
$ws.Range("F45").Value = 1385
$ws.Range("F48").Value = 248
